# The upstream change (commit "Fixed POI packaging and upgraded to POI 3.15")
# only affects how the test-fixture .docx was re-serialized by the Apache POI
# build tooling: every hunk in the diff is a pure XML attribute / namespace
# reordering (attributes alphabetized) on <w:document>, <w:sectPr>, the
# <w:rFonts>/<w:lang> doc defaults, and the <w:latentStyles>/<w:lsdException>/
# <w:style> entries in styles.xml. Comparing the attribute name/value sets
# before and after confirms there is no actual content, formatting, or
# structural change to the document itself.
#
# Since there is nothing to change in the Word object model (no text,
# styles, or layout differ), this script intentionally performs no
# mutation of the document.
$d = $word.ActiveDocument
